$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.966.17"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "3.344.48"
$ws.Range("E3").Value = "  +1.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.29"
$ws.Range("E5").Value = "  +3.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.06"
$ws.Range("E6").Value = "  -5.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "3.340.33"
$ws.Range("E8").Value = "  +1.50%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.613"
$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.08"
$ws.Range("E11").Value = "  -10.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.31"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").Value = "3.881.78"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.344.18"
$ws.Range("E16").Value = "  +2.10%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.53"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "63.881.05"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.27"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "373.37"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").Value = "  +8.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.77"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.69"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.38"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  -0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.98"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "647.90"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.54"
$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.29"
$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.94"
$ws.Range("E36").Value = "  -3.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.19"
$ws.Range("E38").Value = "  +1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.386"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  +10.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  +6.42%  "

$ws.Range("D44").Value = "2.939.82"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.99"
$ws.Range("E45").Value = "  +3.82%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  +3.19%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0401"
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.12"
$ws.Range("E48").Value = "  +5.31%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.64"
$ws.Range("E49").Value = "  -3.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.125"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.83"
$ws.Range("E51").Value = "  +3.70%  "
